$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 701; existing rows 701.. shift down to 703..
$ws.Rows.Item(701).Insert()
$ws.Rows.Item(701).Insert()

# New row 701 (Primera)
$ws.Cells.Item(701,1).Value = 7
$ws.Cells.Item(701,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(701,3).Value = "Ñuble"
$ws.Cells.Item(701,4).Value = "2023-08-09"
$ws.Cells.Item(701,5).Value = 16
$ws.Cells.Item(701,6).Value = 100114014
$ws.Cells.Item(701,7).Value = "Betarraga"
$ws.Cells.Item(701,8).Value = "Sin especificar"
$ws.Cells.Item(701,9).Value = "Primera"
$ws.Cells.Item(701,10).Value = 200
$ws.Cells.Item(701,11).Value = 1000
$ws.Cells.Item(701,12).Value = 1000
$ws.Cells.Item(701,13).Value = 1000
$ws.Cells.Item(701,14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(701,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(701,16).Value = 200
$ws.Cells.Item(701,17).Value = 5
$ws.Cells.Item(701,18).Value = "Hortaliza"

# New row 702 (Segunda)
$ws.Cells.Item(702,1).Value = 7
$ws.Cells.Item(702,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(702,3).Value = "Ñuble"
$ws.Cells.Item(702,4).Value = "2023-08-09"
$ws.Cells.Item(702,5).Value = 16
$ws.Cells.Item(702,6).Value = 100114014
$ws.Cells.Item(702,7).Value = "Betarraga"
$ws.Cells.Item(702,8).Value = "Sin especificar"
$ws.Cells.Item(702,9).Value = "Segunda"
$ws.Cells.Item(702,10).Value = 250
$ws.Cells.Item(702,11).Value = 700
$ws.Cells.Item(702,12).Value = 700
$ws.Cells.Item(702,13).Value = 700
$ws.Cells.Item(702,14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(702,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(702,16).Value = 140
$ws.Cells.Item(702,17).Value = 5
$ws.Cells.Item(702,18).Value = "Hortaliza"
